$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.92"
$ws.Range("E2").Value = "'-0.99%"
$ws.Range("D3").Value = "'31.05"
$ws.Range("E3").Value = "'1.05%"
$ws.Range("D4").Value = "'4.921"
$ws.Range("E4").Value = "'-0.61%"
$ws.Range("D5").Value = "'0.07315"
$ws.Range("E5").Value = "'1.35%"
$ws.Range("D6").Value = "'2.329"
$ws.Range("E6").Value = "'25.00%"
$ws.Range("D7").Value = "'7.741"
$ws.Range("E7").Value = "'0.63%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9033"
$ws.Range("E8").Value = "'0.79%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.09115"
$ws.Range("E9").Value = "'18.39%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1688"
$ws.Range("E10").Value = "'1.66%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08185"
$ws.Range("E11").Value = "'2.33%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03119"
$ws.Range("E12").Value = "'2.77%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09927"
$ws.Range("E13").Value = "'-0.74%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001498"
$ws.Range("E14").Value = "'0.43%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005716"
$ws.Range("E15").Value = "'-0.23%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.494"
$ws.Range("E16").Value = "'0.75%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.720"
$ws.Range("E17").Value = "'-1.43%"
$ws.Range("E19").Value = "'0.36%"
$ws.Range("D20").Value = "'0.1305"
$ws.Range("E20").Value = "'0.56%"
$ws.Range("D21").Value = "'4.215"
$ws.Range("E21").Value = "'4.05%"
$ws.Range("D22").Value = "'0.2100"
$ws.Range("E22").Value = "'-9.52%"
$ws.Range("D23").Value = "'0.04500"
$ws.Range("E23").Value = "'-0.31%"
$ws.Range("D24").Value = "'0.001208"
$ws.Range("E24").Value = "'-0.51%"
$ws.Range("D25").Value = "'0.004162"
$ws.Range("E25").Value = "'-10.30%"
$ws.Range("E26").Value = "'3.93%"
$ws.Range("E39").Value = "'-0.51%"
$ws.Range("D40").Value = "'0.04441"
$ws.Range("D41").Value = "'0.007328"
$ws.Range("E41").Value = "'-0.79%"
$ws.Range("D42").Value = "'0.009514"
$ws.Range("E42").Value = "'-5.19%"
$ws.Range("D43").Value = "'0.1326"
$ws.Range("E43").Value = "'1.74%"
$ws.Range("D44").Value = "'0.002220"
$ws.Range("E44").Value = "'7.17%"
$ws.Range("D45").Value = "'0.008942"
$ws.Range("D46").Value = "'0.00006119"
$ws.Range("E46").Value = "'2.59%"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("D48").Value = "'2.443"
$ws.Range("E48").Value = "'5.98%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("E51").Value = "'-0.04%"
